$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: Swap F:V content between paired/rotated rows (index and date columns A,E unchanged) ---
# Row 4
$ws.Range("F4").Value = 'Unionistas'
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 'Sestao'
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1.67
$ws.Range("K4").Value = '24/08/2023 09:13'
$ws.Range("L4").Value = 1.96
$ws.Range("M4").Value = '26/08/2023 20:07'
$ws.Range("N4").Value = 3.36
$ws.Range("O4").Value = '24/08/2023 09:13'
$ws.Range("P4").Value = 3.24
$ws.Range("Q4").Value = '26/08/2023 20:07'
$ws.Range("R4").Value = 4.85
$ws.Range("S4").Value = '24/08/2023 09:13'
$ws.Range("T4").Value = 4.18
$ws.Range("U4").Value = '26/08/2023 20:07'
$ws.Range("V4").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/unionistas-de-salamanca-sestao/pEPYuf4C/'

# Row 5
$ws.Range("F5").Value = 'Gimnastic'
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 'Arenteiro'
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1.67
$ws.Range("K5").Value = '24/08/2023 09:13'
$ws.Range("L5").Value = 1.73
$ws.Range("M5").Value = '26/08/2023 21:29'
$ws.Range("N5").Value = 3.34
$ws.Range("O5").Value = '24/08/2023 09:13'
$ws.Range("P5").Value = 3.43
$ws.Range("Q5").Value = '26/08/2023 21:29'
$ws.Range("R5").Value = 4.91
$ws.Range("S5").Value = '24/08/2023 09:13'
$ws.Range("T5").Value = 5.32
$ws.Range("U5").Value = '26/08/2023 21:29'
$ws.Range("V5").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/gimnastic-de-tarragona-arenteiro/4IheceZn/'

# Row 9
$ws.Range("F9").Value = 'SD Logrones'
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 'Barcelona B'
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 2.3
$ws.Range("K9").Value = '24/08/2023 09:13'
$ws.Range("L9").Value = 2.96
$ws.Range("M9").Value = '27/08/2023 18:29'
$ws.Range("N9").Value = 3.08
$ws.Range("O9").Value = '24/08/2023 09:13'
$ws.Range("P9").Value = 3.22
$ws.Range("Q9").Value = '27/08/2023 19:33'
$ws.Range("R9").Value = 2.93
$ws.Range("S9").Value = '24/08/2023 09:13'
$ws.Range("T9").Value = 2.36
$ws.Range("U9").Value = '27/08/2023 18:29'
$ws.Range("V9").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/sd-logrones-barcelona/WUNtvYYO/'

# Row 10
$ws.Range("F10").Value = 'Ponferradina'
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 'Celta Vigo B'
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 2.08
$ws.Range("K10").Value = '27/08/2023 16:43'
$ws.Range("L10").Value = 2.08
$ws.Range("M10").Value = '27/08/2023 20:56'
$ws.Range("N10").Value = 3.37
$ws.Range("O10").Value = '27/08/2023 16:43'
$ws.Range("P10").Value = 3.43
$ws.Range("Q10").Value = '27/08/2023 21:01'
$ws.Range("R10").Value = 3.48
$ws.Range("S10").Value = '27/08/2023 16:43'
$ws.Range("T10").Value = 3.5
$ws.Range("U10").Value = '27/08/2023 20:56'
$ws.Range("V10").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/ponferradina-celta-vigo/r9iadFlg/'

# Row 11
$ws.Range("F11").Value = 'Fuenlabrada'
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 'Leonesa'
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 2.44
$ws.Range("K11").Value = '26/08/2023 05:12'
$ws.Range("L11").Value = 3.23
$ws.Range("M11").Value = '27/08/2023 11:02'
$ws.Range("N11").Value = 2.98
$ws.Range("O11").Value = '26/08/2023 05:12'
$ws.Range("P11").Value = 2.97
$ws.Range("Q11").Value = '27/08/2023 19:32'
$ws.Range("R11").Value = 2.81
$ws.Range("S11").Value = '26/08/2023 05:12'
$ws.Range("T11").Value = 2.41
$ws.Range("U11").Value = '27/08/2023 11:02'
$ws.Range("V11").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/cf-fuenlabrada-leonesa/6NQUtzk6/'

# Row 14
$ws.Range("F14").Value = 'Rayo Majadahonda'
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 'Gimnastic'
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 2.46
$ws.Range("K14").Value = '31/08/2023 09:13'
$ws.Range("L14").Value = 3.38
$ws.Range("M14").Value = '02/09/2023 21:29'
$ws.Range("N14").Value = 2.93
$ws.Range("O14").Value = '31/08/2023 09:13'
$ws.Range("P14").Value = 3.12
$ws.Range("Q14").Value = '02/09/2023 21:29'
$ws.Range("R14").Value = 2.84
$ws.Range("S14").Value = '31/08/2023 09:13'
$ws.Range("T14").Value = 2.27
$ws.Range("U14").Value = '02/09/2023 21:29'
$ws.Range("V14").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/rayo-majadahonda-gimnastic-de-tarragona/E7hKii3O/'

# Row 15
$ws.Range("F15").Value = 'Lugo'
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 'Dep. La Coruna'
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 2.96
$ws.Range("K15").Value = '31/08/2023 09:13'
$ws.Range("L15").Value = 3.26
$ws.Range("M15").Value = '02/09/2023 21:20'
$ws.Range("N15").Value = 3.04
$ws.Range("O15").Value = '31/08/2023 09:13'
$ws.Range("P15").Value = 3.11
$ws.Range("Q15").Value = '02/09/2023 21:20'
$ws.Range("R15").Value = 2.35
$ws.Range("S15").Value = '31/08/2023 09:13'
$ws.Range("T15").Value = 2.33
$ws.Range("U15").Value = '02/09/2023 21:20'
$ws.Range("V15").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/lugo-dep-la-coruna/hbiOjBIU/'

# Row 16
$ws.Range("F16").Value = 'Sestao'
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 'Ponferradina'
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 2.67
$ws.Range("K16").Value = '31/08/2023 09:13'
$ws.Range("L16").Value = 2.78
$ws.Range("M16").Value = '03/09/2023 18:36'
$ws.Range("N16").Value = 3.1
$ws.Range("O16").Value = '31/08/2023 09:13'
$ws.Range("P16").Value = 3.06
$ws.Range("Q16").Value = '03/09/2023 18:36'
$ws.Range("R16").Value = 2.48
$ws.Range("S16").Value = '31/08/2023 09:13'
$ws.Range("T16").Value = 2.7
$ws.Range("U16").Value = '03/09/2023 18:36'
$ws.Range("V16").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/sestao-ponferradina/QoMxAYIh/'

# Row 18
$ws.Range("F18").Value = 'Arenteiro'
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 'SD Logrones'
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 2.25
$ws.Range("K18").Value = '31/08/2023 09:13'
$ws.Range("L18").Value = 2.32
$ws.Range("M18").Value = '03/09/2023 15:25'
$ws.Range("N18").Value = 2.98
$ws.Range("O18").Value = '31/08/2023 09:13'
$ws.Range("P18").Value = 2.9
$ws.Range("Q18").Value = '03/09/2023 18:34'
$ws.Range("R18").Value = 3.11
$ws.Range("S18").Value = '31/08/2023 09:13'
$ws.Range("T18").Value = 3.5
$ws.Range("U18").Value = '03/09/2023 15:25'
$ws.Range("V18").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/arenteiro-sd-logrones/YcpCgDYB/'

# Row 19
$ws.Range("F19").Value = 'Teruel'
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 'R. Sociedad B'
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2.46
$ws.Range("K19").Value = '31/08/2023 09:13'
$ws.Range("L19").Value = 2.66
$ws.Range("M19").Value = '03/09/2023 21:20'
$ws.Range("N19").Value = 2.93
$ws.Range("O19").Value = '31/08/2023 09:13'
$ws.Range("P19").Value = 3.05
$ws.Range("Q19").Value = '03/09/2023 21:20'
$ws.Range("R19").Value = 2.84
$ws.Range("S19").Value = '31/08/2023 09:13'
$ws.Range("T19").Value = 2.83
$ws.Range("U19").Value = '03/09/2023 21:20'
$ws.Range("V19").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/teruel-r-sociedad/vkQt9hYb/'

# Row 20
$ws.Range("F20").Value = 'Sabadell'
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 'Osasuna B'
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 2.46
$ws.Range("K20").Value = '31/08/2023 09:13'
$ws.Range("L20").Value = 2.22
$ws.Range("M20").Value = '03/09/2023 21:28'
$ws.Range("N20").Value = 2.93
$ws.Range("O20").Value = '31/08/2023 09:13'
$ws.Range("P20").Value = 3.12
$ws.Range("Q20").Value = '03/09/2023 21:21'
$ws.Range("R20").Value = 2.84
$ws.Range("S20").Value = '31/08/2023 09:13'
$ws.Range("T20").Value = 3.5
$ws.Range("U20").Value = '03/09/2023 21:28'
$ws.Range("V20").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/sabadell-osasuna/4WXTBflt/'

# Row 21
$ws.Range("F21").Value = 'Cornella'
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 'Fuenlabrada'
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 2.25
$ws.Range("K21").Value = '31/08/2023 09:13'
$ws.Range("L21").Value = 2.28
$ws.Range("M21").Value = '03/09/2023 21:23'
$ws.Range("N21").Value = 2.98
$ws.Range("O21").Value = '31/08/2023 09:13'
$ws.Range("P21").Value = 2.96
$ws.Range("Q21").Value = '03/09/2023 21:23'
$ws.Range("R21").Value = 3.11
$ws.Range("S21").Value = '31/08/2023 09:13'
$ws.Range("T21").Value = 3.55
$ws.Range("U21").Value = '03/09/2023 21:23'
$ws.Range("V21").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/cornella-cf-fuenlabrada/vyo8fgJ5/'

# Row 22
$ws.Range("F22").Value = 'Tarazona'
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 'Leonesa'
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2.91
$ws.Range("K22").Value = '07/09/2023 09:13'
$ws.Range("L22").Value = 4.1
$ws.Range("M22").Value = '09/09/2023 15:54'
$ws.Range("N22").Value = 2.99
$ws.Range("O22").Value = '07/09/2023 09:13'
$ws.Range("P22").Value = 3.02
$ws.Range("Q22").Value = '09/09/2023 17:32'
$ws.Range("R22").Value = 2.36
$ws.Range("S22").Value = '07/09/2023 09:13'
$ws.Range("T22").Value = 2.05
$ws.Range("U22").Value = '09/09/2023 15:54'
$ws.Range("V22").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/tarazona-leonesa/hxbE5vq4/'

# Row 23
$ws.Range("F23").Value = 'Ponferradina'
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 'Sabadell'
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 1.91
$ws.Range("K23").Value = '07/09/2023 09:13'
$ws.Range("L23").Value = 1.96
$ws.Range("M23").Value = '09/09/2023 16:10'
$ws.Range("N23").Value = 3.2
$ws.Range("O23").Value = '07/09/2023 09:13'
$ws.Range("P23").Value = 2.94
$ws.Range("Q23").Value = '09/09/2023 17:32'
$ws.Range("R23").Value = 3.76
$ws.Range("S23").Value = '07/09/2023 09:13'
$ws.Range("T23").Value = 4.69
$ws.Range("U23").Value = '09/09/2023 16:10'
$ws.Range("V23").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/ponferradina-sabadell/W4vj5AXN/'

# Row 26
$ws.Range("F26").Value = 'Fuenlabrada'
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 'Lugo'
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 2.25
$ws.Range("K26").Value = '07/09/2023 09:13'
$ws.Range("L26").Value = 2.63
$ws.Range("M26").Value = '10/09/2023 11:57'
$ws.Range("N26").Value = 3.02
$ws.Range("O26").Value = '07/09/2023 09:13'
$ws.Range("P26").Value = 3
$ws.Range("Q26").Value = '10/09/2023 11:57'
$ws.Range("R26").Value = 3.06
$ws.Range("S26").Value = '07/09/2023 09:13'
$ws.Range("T26").Value = 2.91
$ws.Range("U26").Value = '10/09/2023 11:57'
$ws.Range("V26").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/cf-fuenlabrada-lugo/hEun6jIH/'

# Row 27
$ws.Range("F27").Value = 'R. Sociedad B'
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 'Unionistas'
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 1.95
$ws.Range("K27").Value = '07/09/2023 09:13'
$ws.Range("L27").Value = 2.44
$ws.Range("M27").Value = '10/09/2023 11:58'
$ws.Range("N27").Value = 3.07
$ws.Range("O27").Value = '07/09/2023 09:13'
$ws.Range("P27").Value = 3.07
$ws.Range("Q27").Value = '10/09/2023 11:58'
$ws.Range("R27").Value = 3.79
$ws.Range("S27").Value = '07/09/2023 09:13'
$ws.Range("T27").Value = 3.11
$ws.Range("U27").Value = '10/09/2023 11:58'
$ws.Range("V27").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/r-sociedad-unionistas-de-salamanca/2mVLaleo/'

# Row 29
$ws.Range("F29").Value = 'Gimnastic'
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 'Barcelona B'
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2.1
$ws.Range("K29").Value = '07/09/2023 09:13'
$ws.Range("L29").Value = 2.15
$ws.Range("M29").Value = '10/09/2023 19:20'
$ws.Range("N29").Value = 3.11
$ws.Range("O29").Value = '07/09/2023 09:13'
$ws.Range("P29").Value = 3.25
$ws.Range("Q29").Value = '10/09/2023 19:20'
$ws.Range("R29").Value = 3.3
$ws.Range("S29").Value = '07/09/2023 09:13'
$ws.Range("T29").Value = 3.52
$ws.Range("U29").Value = '10/09/2023 19:20'
$ws.Range("V29").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/gimnastic-de-tarragona-barcelona/6Le28tao/'

# Row 30
$ws.Range("F30").Value = 'Dep. La Coruna'
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 'Teruel'
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1.39
$ws.Range("K30").Value = '07/09/2023 09:13'
$ws.Range("L30").Value = 1.29
$ws.Range("M30").Value = '10/09/2023 19:23'
$ws.Range("N30").Value = 4.27
$ws.Range("O30").Value = '07/09/2023 09:13'
$ws.Range("P30").Value = 4.93
$ws.Range("Q30").Value = '10/09/2023 19:28'
$ws.Range("R30").Value = 7.38
$ws.Range("S30").Value = '07/09/2023 09:13'
$ws.Range("T30").Value = 13.15
$ws.Range("U30").Value = '10/09/2023 19:28'
$ws.Range("V30").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/dep-la-coruna-teruel/YBPp8Cm5/'

# Row 31
$ws.Range("F31").Value = 'R. Union'
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 'Rayo Majadahonda'
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1.65
$ws.Range("K31").Value = '07/09/2023 09:13'
$ws.Range("L31").Value = 1.71
$ws.Range("M31").Value = '10/09/2023 19:21'
$ws.Range("N31").Value = 3.49
$ws.Range("O31").Value = '07/09/2023 09:13'
$ws.Range("P31").Value = 3.55
$ws.Range("Q31").Value = '10/09/2023 19:21'
$ws.Range("R31").Value = 4.98
$ws.Range("S31").Value = '07/09/2023 09:13'
$ws.Range("T31").Value = 5.23
$ws.Range("U31").Value = '10/09/2023 19:21'
$ws.Range("V31").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/real-union-rayo-majadahonda/EoaA6KTc/'

# Row 37
$ws.Range("F37").Value = 'Leonesa'
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 'Osasuna B'
$ws.Range("I37").Value = 2
$ws.Range("J37").Value = 1.94
$ws.Range("K37").Value = '14/09/2023 15:42'
$ws.Range("L37").Value = 1.93
$ws.Range("M37").Value = '17/09/2023 11:51'
$ws.Range("N37").Value = 3.23
$ws.Range("O37").Value = '14/09/2023 15:42'
$ws.Range("P37").Value = 3.4
$ws.Range("Q37").Value = '17/09/2023 11:51'
$ws.Range("R37").Value = 3.61
$ws.Range("S37").Value = '14/09/2023 15:42'
$ws.Range("T37").Value = 4.08
$ws.Range("U37").Value = '17/09/2023 11:51'
$ws.Range("V37").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/leonesa-osasuna/za8apUCe/'

# Row 38
$ws.Range("F38").Value = 'Cornella'
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 'Tarazona'
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 1.74
$ws.Range("K38").Value = '14/09/2023 15:42'
$ws.Range("L38").Value = 1.85
$ws.Range("M38").Value = '16/09/2023 17:35'
$ws.Range("N38").Value = 3.26
$ws.Range("O38").Value = '14/09/2023 15:42'
$ws.Range("P38").Value = 3.3
$ws.Range("Q38").Value = '17/09/2023 10:05'
$ws.Range("R38").Value = 4.57
$ws.Range("S38").Value = '14/09/2023 15:42'
$ws.Range("T38").Value = 4.64
$ws.Range("U38").Value = '17/09/2023 10:29'
$ws.Range("V38").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/cornella-tarazona/nyiN3IDG/'

# Row 70
$ws.Range("F70").Value = 'R. Sociedad B'
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 'Barcelona B'
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2.35
$ws.Range("K70").Value = '05/10/2023 08:13'
$ws.Range("L70").Value = 1.82
$ws.Range("M70").Value = '08/10/2023 19:50'
$ws.Range("N70").Value = 2.97
$ws.Range("O70").Value = '05/10/2023 08:13'
$ws.Range("P70").Value = 3.39
$ws.Range("Q70").Value = '08/10/2023 19:59'
$ws.Range("R70").Value = 2.94
$ws.Range("S70").Value = '05/10/2023 08:13'
$ws.Range("T70").Value = 4.69
$ws.Range("U70").Value = '08/10/2023 19:59'
$ws.Range("V70").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/r-sociedad-barcelona/t0sfN3A6/'

# Row 71
$ws.Range("F71").Value = 'Osasuna B'
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 'Lugo'
$ws.Range("I71").Value = 2
$ws.Range("J71").Value = 2.15
$ws.Range("K71").Value = '05/10/2023 08:13'
$ws.Range("L71").Value = 2.33
$ws.Range("M71").Value = '08/10/2023 19:54'
$ws.Range("N71").Value = 3.01
$ws.Range("O71").Value = '05/10/2023 08:13'
$ws.Range("P71").Value = 3.07
$ws.Range("Q71").Value = '08/10/2023 19:54'
$ws.Range("R71").Value = 3.29
$ws.Range("S71").Value = '05/10/2023 08:13'
$ws.Range("T71").Value = 3.3
$ws.Range("U71").Value = '08/10/2023 19:54'
$ws.Range("V71").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/osasuna-lugo/lQp2LsuJ/'

# Row 86
$ws.Range("F86").Value = 'Osasuna B'
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 'Celta Vigo B'
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 2.38
$ws.Range("K86").Value = '19/10/2023 08:13'
$ws.Range("L86").Value = 2.41
$ws.Range("M86").Value = '22/10/2023 11:52'
$ws.Range("N86").Value = 3.12
$ws.Range("O86").Value = '19/10/2023 08:13'
$ws.Range("P86").Value = 3.06
$ws.Range("Q86").Value = '22/10/2023 11:52'
$ws.Range("R86").Value = 2.85
$ws.Range("S86").Value = '19/10/2023 08:13'
$ws.Range("T86").Value = 3.17
$ws.Range("U86").Value = '22/10/2023 11:52'
$ws.Range("V86").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/osasuna-celta-vigo/lrvwhFFc/'

# Row 87
$ws.Range("F87").Value = 'SD Logrones'
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 'Cornella'
$ws.Range("I87").Value = 2
$ws.Range("J87").Value = 2.23
$ws.Range("K87").Value = '19/10/2023 08:13'
$ws.Range("L87").Value = 2.55
$ws.Range("M87").Value = '22/10/2023 11:57'
$ws.Range("N87").Value = 2.96
$ws.Range("O87").Value = '19/10/2023 08:13'
$ws.Range("P87").Value = 2.87
$ws.Range("Q87").Value = '22/10/2023 11:57'
$ws.Range("R87").Value = 3.18
$ws.Range("S87").Value = '19/10/2023 08:13'
$ws.Range("T87").Value = 3.16
$ws.Range("U87").Value = '22/10/2023 11:57'
$ws.Range("V87").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/sd-logrones-cornella/vXAKZXaT/'

# Row 89
$ws.Range("F89").Value = 'Gimnastic'
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 'R. Sociedad B'
$ws.Range("I89").Value = 2
$ws.Range("J89").Value = 1.9
$ws.Range("K89").Value = '19/10/2023 08:13'
$ws.Range("L89").Value = 1.97
$ws.Range("M89").Value = '22/10/2023 16:10'
$ws.Range("N89").Value = 3.11
$ws.Range("O89").Value = '19/10/2023 08:13'
$ws.Range("P89").Value = 3.16
$ws.Range("Q89").Value = '22/10/2023 16:10'
$ws.Range("R89").Value = 3.93
$ws.Range("S89").Value = '19/10/2023 08:13'
$ws.Range("T89").Value = 4.26
$ws.Range("U89").Value = '22/10/2023 16:10'
$ws.Range("V89").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/gimnastic-de-tarragona-r-sociedad/4EI7xzGA/'

# Row 90
$ws.Range("F90").Value = 'Rayo Majadahonda'
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 'Lugo'
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 2.65
$ws.Range("K90").Value = '19/10/2023 08:13'
$ws.Range("L90").Value = 3
$ws.Range("M90").Value = '22/10/2023 16:53'
$ws.Range("N90").Value = 2.88
$ws.Range("O90").Value = '19/10/2023 08:13'
$ws.Range("P90").Value = 3
$ws.Range("Q90").Value = '22/10/2023 16:53'
$ws.Range("R90").Value = 2.65
$ws.Range("S90").Value = '19/10/2023 08:13'
$ws.Range("T90").Value = 2.56
$ws.Range("U90").Value = '22/10/2023 16:53'
$ws.Range("V90").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/rayo-majadahonda-lugo/8n1jUBNj/'

# Row 92
$ws.Range("F92").Value = 'Celta Vigo B'
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 'Sabadell'
$ws.Range("I92").Value = 2
$ws.Range("J92").Value = 1.82
$ws.Range("K92").Value = '26/10/2023 23:12'
$ws.Range("L92").Value = 1.88
$ws.Range("M92").Value = '28/10/2023 13:35'
$ws.Range("N92").Value = 3.37
$ws.Range("O92").Value = '26/10/2023 23:12'
$ws.Range("P92").Value = 3.5
$ws.Range("Q92").Value = '28/10/2023 15:01'
$ws.Range("R92").Value = 3.9
$ws.Range("S92").Value = '26/10/2023 23:12'
$ws.Range("T92").Value = 4.18
$ws.Range("U92").Value = '28/10/2023 15:28'
$ws.Range("V92").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/celta-vigo-sabadell/fJiQtk7d/'

# Row 93
$ws.Range("F93").Value = 'Arenteiro'
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 'Unionistas'
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2.34
$ws.Range("K93").Value = '26/10/2023 23:12'
$ws.Range("L93").Value = 2.92
$ws.Range("M93").Value = '28/10/2023 15:39'
$ws.Range("N93").Value = 2.91
$ws.Range("O93").Value = '26/10/2023 23:12'
$ws.Range("P93").Value = 2.9
$ws.Range("Q93").Value = '28/10/2023 15:39'
$ws.Range("R93").Value = 3.02
$ws.Range("S93").Value = '26/10/2023 23:12'
$ws.Range("T93").Value = 2.7
$ws.Range("U93").Value = '28/10/2023 15:39'
$ws.Range("V93").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/arenteiro-unionistas-de-salamanca/GShMsVij/'

# Row 96
$ws.Range("F96").Value = 'Tarazona'
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 'Sestao'
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2.13
$ws.Range("K96").Value = '26/10/2023 23:12'
$ws.Range("L96").Value = 2.45
$ws.Range("M96").Value = '29/10/2023 15:56'
$ws.Range("N96").Value = 3.05
$ws.Range("O96").Value = '26/10/2023 23:12'
$ws.Range("P96").Value = 2.7
$ws.Range("Q96").Value = '29/10/2023 15:56'
$ws.Range("R96").Value = 3.29
$ws.Range("S96").Value = '26/10/2023 23:12'
$ws.Range("T96").Value = 3.59
$ws.Range("U96").Value = '29/10/2023 15:56'
$ws.Range("V96").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/tarazona-sestao/ruYK2WTF/'

# Row 97
$ws.Range("F97").Value = 'Cornella'
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 'Teruel'
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 2.14
$ws.Range("K97").Value = '26/10/2023 23:12'
$ws.Range("L97").Value = 2.28
$ws.Range("M97").Value = '29/10/2023 15:51'
$ws.Range("N97").Value = 2.96
$ws.Range("O97").Value = '26/10/2023 23:12'
$ws.Range("P97").Value = 2.75
$ws.Range("Q97").Value = '29/10/2023 15:51'
$ws.Range("R97").Value = 3.36
$ws.Range("S97").Value = '26/10/2023 23:12'
$ws.Range("T97").Value = 3.9
$ws.Range("U97").Value = '29/10/2023 15:51'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/cornella-teruel/hzyflXEM/'

# --- Part 2: Append new rows 102-104 (copy style from row 101, then set values) ---
# Row 102
$ws.Range("A101:V101").Copy($ws.Range("A102:V102"))
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 'spain'
$ws.Range("C102").Value = 'primera-rfef-group-1'
$ws.Range("D102").Value = '2023-2024'
$ws.Range("E102").Value = 45234.70833333334
$ws.Range("F102").Value = 'Rayo Majadahonda'
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 'Barcelona B'
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 3.64
$ws.Range("K102").Value = '02/11/2023 08:13'
$ws.Range("L102").Value = 2.94
$ws.Range("M102").Value = '04/11/2023 16:57'
$ws.Range("N102").Value = 3.2
$ws.Range("O102").Value = '02/11/2023 08:13'
$ws.Range("P102").Value = 3
$ws.Range("Q102").Value = '04/11/2023 16:57'
$ws.Range("R102").Value = 1.98
$ws.Range("S102").Value = '02/11/2023 08:13'
$ws.Range("T102").Value = 2.61
$ws.Range("U102").Value = '04/11/2023 16:57'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/rayo-majadahonda-barcelona/dnEedU6q/'

# Row 103
$ws.Range("A101:V101").Copy($ws.Range("A103:V103"))
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = 'spain'
$ws.Range("C103").Value = 'primera-rfef-group-1'
$ws.Range("D103").Value = '2023-2024'
$ws.Range("E103").Value = 45234.79166666666
$ws.Range("F103").Value = 'Arenteiro'
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 'Fuenlabrada'
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = 2.25
$ws.Range("K103").Value = '02/11/2023 08:13'
$ws.Range("L103").Value = 2.33
$ws.Range("M103").Value = '04/11/2023 18:51'
$ws.Range("N103").Value = 2.98
$ws.Range("O103").Value = '02/11/2023 08:13'
$ws.Range("P103").Value = 2.89
$ws.Range("Q103").Value = '04/11/2023 18:51'
$ws.Range("R103").Value = 3.11
$ws.Range("S103").Value = '02/11/2023 08:13'
$ws.Range("T103").Value = 3.55
$ws.Range("U103").Value = '04/11/2023 18:51'
$ws.Range("V103").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/arenteiro-cf-fuenlabrada/CExP1jqM/'

# Row 104
$ws.Range("A101:V101").Copy($ws.Range("A104:V104"))
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = 'spain'
$ws.Range("C104").Value = 'primera-rfef-group-1'
$ws.Range("D104").Value = '2023-2024'
$ws.Range("E104").Value = 45234.83333333334
$ws.Range("F104").Value = 'Leonesa'
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 'Cornella'
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 1.79
$ws.Range("K104").Value = '02/11/2023 08:13'
$ws.Range("L104").Value = 1.7
$ws.Range("M104").Value = '04/11/2023 19:51'
$ws.Range("N104").Value = 3.21
$ws.Range("O104").Value = '02/11/2023 08:13'
$ws.Range("P104").Value = 3.5
$ws.Range("Q104").Value = '04/11/2023 19:55'
$ws.Range("R104").Value = 4.34
$ws.Range("S104").Value = '02/11/2023 08:13'
$ws.Range("T104").Value = 5.4
$ws.Range("U104").Value = '04/11/2023 19:51'
$ws.Range("V104").Value = 'https://www.betexplorer.com/football/spain/primera-rfef-group-1/leonesa-cornella/8MHaelMk/'

